$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dwellings_buildings")

# Row 2 (H:1 - Mining and quarrying -> All other industry)
$ws.Range("B2").Value = "All other industry"
$ws.Range("C2").Value = 220
$ws.Range("D2").Value = 317.6267100977199

# Row 3 (H:1 - Transport; storage and communication -> Warehouses and storage)
$ws.Range("B3").Value = "Warehouses and storage"
$ws.Range("C3").Value = 220
$ws.Range("D3").Value = 245.5050488599348

# Row 4 (H:1 - Manufacturing -> Manufacturing and light industry)
$ws.Range("B4").Value = "Manufacturing and light industry"
$ws.Range("C4").Value = 220
$ws.Range("D4").Value = 295.0431596091205

# Row 5 (H:2 - Mining and quarrying -> All other industry)
$ws.Range("B5").Value = "All other industry"
$ws.Range("C5").Value = 480
$ws.Range("D5").Value = 317.6267100977199

# Row 6 (H:2 - Transport; storage and communication -> Warehouses and storage)
$ws.Range("B6").Value = "Warehouses and storage"
$ws.Range("C6").Value = 480
$ws.Range("D6").Value = 245.5050488599348

# Row 7 (H:2 - Manufacturing -> Manufacturing and light industry)
$ws.Range("B7").Value = "Manufacturing and light industry"
$ws.Range("C7").Value = 480
$ws.Range("D7").Value = 295.0431596091205
